# Generate Report for Handback
# Updates the localization-status report: the zh-cn file is now fully
# handed back (in sync with en-US) and its handback error has cleared;
# the de-de file's handback timestamp advanced too.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status text changed everywhere it was shown: Overview (both language
# columns), and the Status column on each per-language sheet.
$newStatus = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# zh-cn: handback completed successfully -> no more error detail, and the
# "Latest Handback DateTime" advances to the new handback time.
$wsZhCn.Range("K2").Value = "2016-09-04 20:56:17"
$wsZhCn.Range("P2").Value = ""

# de-de: same story - handback datetime advances and the stale error
# detail is cleared.
$wsDeDe.Range("K2").Value = "2016-09-04 20:56:24"
$wsDeDe.Range("P2").Value = ""

# Column widths grew to fit the new (longer) Status text, and the now
# mostly-empty "Error Detail" column shrank back down. (Input values are
# pre-snapped to this host's 1/6-character ColumnWidth grid so the saved
# OOXML width lands as close as possible to the target width.)
$wsOverview.Columns.Item(5).ColumnWidth = 29.1666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.1666666666667

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsZhCn.Columns.Item(16).ColumnWidth = 12.8333333333333

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsDeDe.Columns.Item(16).ColumnWidth = 12.8333333333333
